# msg_8706.xlsx: support for new msgId 8717
# Appends 31 new telemetry rows (rows 563-593) captured 2022-03-22,
# all sharing the same dtuSn ("8: 255") and offset ("10F872226797") as
# the prior block, continuing columns C (period=3600), D (unix ts),
# E (rssi), G (val), H (val7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ R=563; A="2022-03-22 13:10:01"; C=3600; D=1647951003; E=56; G=180; H=120 },
    @{ R=564; A="2022-03-22 20:21:31"; C=3600; D=1647976895; E=46; G=-59562; H=165 },
    @{ R=565; A="2022-03-22 20:21:52"; C=3600; D=1647976916; E=46; G=238; H=165 },
    @{ R=566; A="2022-03-22 20:22:52"; C=3600; D=1647976976; E=46; G=473; H=165 },
    @{ R=567; A="2022-03-22 20:23:53"; C=3600; D=1647977037; E=46; G=301; H=165 },
    @{ R=568; A="2022-03-22 20:24:53"; C=3600; D=1647977097; E=46; G=206; H=165 },
    @{ R=569; A="2022-03-22 20:25:54"; C=3600; D=1647977158; E=46; G=222; H=165 },
    @{ R=570; A="2022-03-22 20:26:53"; C=3600; D=1647977217; E=46; G=193; H=165 },
    @{ R=571; A="2022-03-22 20:27:53"; C=3600; D=1647977277; E=46; G=199; H=165 },
    @{ R=572; A="2022-03-22 20:28:54"; C=3600; D=1647977338; E=46; G=220; H=165 },
    @{ R=573; A="2022-03-22 20:29:53"; C=3600; D=1647977397; E=46; G=195; H=165 },
    @{ R=574; A="2022-03-22 20:30:54"; C=3600; D=1647977458; E=46; G=245; H=165 },
    @{ R=575; A="2022-03-22 20:31:55"; C=3600; D=1647977519; E=46; G=214; H=165 },
    @{ R=576; A="2022-03-22 20:32:54"; C=3600; D=1647977578; E=46; G=228; H=165 },
    @{ R=577; A="2022-03-22 20:33:54"; C=3600; D=1647977638; E=46; G=213; H=165 },
    @{ R=578; A="2022-03-22 20:34:54"; C=3600; D=1647977698; E=46; G=197; H=165 },
    @{ R=579; A="2022-03-22 20:35:54"; C=3600; D=1647977758; E=46; G=254; H=165 },
    @{ R=580; A="2022-03-22 20:36:54"; C=3600; D=1647977818; E=46; G=260; H=165 },
    @{ R=581; A="2022-03-22 20:37:54"; C=3600; D=1647977878; E=46; G=223; H=165 },
    @{ R=582; A="2022-03-22 20:38:54"; C=3600; D=1647977938; E=46; G=296; H=165 },
    @{ R=583; A="2022-03-22 20:39:54"; C=3600; D=1647977998; E=46; G=291; H=165 },
    @{ R=584; A="2022-03-22 20:40:54"; C=3600; D=1647978058; E=46; G=241; H=165 },
    @{ R=585; A="2022-03-22 20:41:54"; C=3600; D=1647978118; E=46; G=250; H=165 },
    @{ R=586; A="2022-03-22 20:42:55"; C=3600; D=1647978179; E=46; G=218; H=165 },
    @{ R=587; A="2022-03-22 20:43:54"; C=3600; D=1647978238; E=46; G=201; H=165 },
    @{ R=588; A="2022-03-22 20:44:54"; C=3600; D=1647978298; E=46; G=197; H=165 },
    @{ R=589; A="2022-03-22 20:45:55"; C=3600; D=1647978359; E=46; G=234; H=165 },
    @{ R=590; A="2022-03-22 20:46:55"; C=3600; D=1647978419; E=46; G=227; H=165 },
    @{ R=591; A="2022-03-22 20:47:55"; C=3600; D=1647978479; E=46; G=205; H=165 },
    @{ R=592; A="2022-03-22 20:48:55"; C=3600; D=1647978539; E=46; G=212; H=165 },
    @{ R=593; A="2022-03-22 20:49:55"; C=3600; D=1647978599; E=46; G=205; H=165 }
)

foreach ($row in $newRows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.A          # datetime (col A)
    $ws.Cells.Item($r, 2).Value = "8: 255`n"      # dtuSn (col B) - matches existing shared string w/ trailing newline
    $ws.Cells.Item($r, 3).Value = $row.C          # offset/period (col C)
    $ws.Cells.Item($r, 4).Value = $row.D          # unix timestamp (col D)
    $ws.Cells.Item($r, 5).Value = $row.E          # val3 (col E)
    $ws.Cells.Item($r, 6).Value = "10F872226797"  # val6 (col F)
    $ws.Cells.Item($r, 7).Value = $row.G          # val7 (col G)
    $ws.Cells.Item($r, 8).Value = $row.H          # (col H)
}

# The embedded newline in column B makes the host auto-size the row
# (adds ht/customHeight). Re-running AutoFit restores the default,
# un-flagged row height used throughout the rest of the sheet.
$ws.Range("563:593").EntireRow.AutoFit()
